# corrected data cleaning for pre/post/total fixation data
#
# This script:
#  1. Strips the bold/border/center-top formatting from the header row
#     (A1:AL1), matching the removal of the now-unused font/border/cellXf
#     entries from styles.xml.
#  2. Clears the "Unnamed: 0" label from A1 (header becomes blank).
#  3. Clears the entire "declaration" column (O) data values in rows 3-8,
#     since that column is dropped from the cleaned aggregate.
#  4. Updates the numeric aggregate values in the remaining data columns
#     for rows 3-7 to their recomputed (corrected) totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove header formatting (bold font, thin borders, center/top alignment)
$ws.Range("A1:AL1").ClearFormats()

# 2) Clear the old "Unnamed: 0" header text
$ws.Range("A1").Value = ""

# 3) Clear "declaration" (column O) values for the data rows
$ws.Range("O3").Value = ""
$ws.Range("O4").Value = ""
$ws.Range("O5").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("O7").Value = ""
$ws.Range("O8").Value = ""

# 4) Recomputed values for row 3 (Revisit count)
$ws.Range("D3").Value = 12
$ws.Range("F3").Value = 13
$ws.Range("G3").Value = 19
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 57
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 45
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 41
$ws.Range("Q3").Value = 15
$ws.Range("R3").Value = 19
$ws.Range("S3").Value = 11
$ws.Range("V3").Value = 51
$ws.Range("W3").Value = 53
$ws.Range("X3").Value = 54
$ws.Range("AA3").Value = 23
$ws.Range("AE3").Value = 53
$ws.Range("AG3").Value = 6

# Recomputed values for row 4 (Fixation count)
$ws.Range("D4").Value = 18
$ws.Range("F4").Value = 22
$ws.Range("G4").Value = 34
$ws.Range("H4").Value = 248
$ws.Range("I4").Value = 219
$ws.Range("J4").Value = 156
$ws.Range("K4").Value = 107
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 72
$ws.Range("Q4").Value = 16
$ws.Range("R4").Value = 30
$ws.Range("S4").Value = 14
$ws.Range("V4").Value = 243
$ws.Range("W4").Value = 238
$ws.Range("X4").Value = 187
$ws.Range("AA4").Value = 27
$ws.Range("AE4").Value = 220
$ws.Range("AG4").Value = 9

# Recomputed values for row 5 (Dwell time (ms))
$ws.Range("D5").Value = 5672.69
$ws.Range("F5").Value = 9001.870000000001
$ws.Range("G5").Value = 11653.9
$ws.Range("H5").Value = 74744.07000000001
$ws.Range("I5").Value = 65871.62
$ws.Range("J5").Value = 49595.5
$ws.Range("K5").Value = 35136.06
$ws.Range("M5").Value = 17802.68
$ws.Range("N5").Value = 26418.36
$ws.Range("Q5").Value = 8258.809999999999
$ws.Range("R5").Value = 10744.86
$ws.Range("S5").Value = 4571.62
$ws.Range("V5").Value = 72992.35000000001
$ws.Range("W5").Value = 72208.33
$ws.Range("X5").Value = 59106.17
$ws.Range("AA5").Value = 11835.9
$ws.Range("AE5").Value = 90194.37
$ws.Range("AG5").Value = 4204.46

# Recomputed values for row 6 (Dwell time (%))
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0.07000000000000001
$ws.Range("D6").Value = 2.19
$ws.Range("F6").Value = 3.48
$ws.Range("G6").Value = 4.51
$ws.Range("H6").Value = 28.91
$ws.Range("I6").Value = 25.48
$ws.Range("J6").Value = 19.18
$ws.Range("K6").Value = 13.59
$ws.Range("L6").Value = 2.17
$ws.Range("M6").Value = 6.89
$ws.Range("N6").Value = 10.22
$ws.Range("Q6").Value = 3.19
$ws.Range("R6").Value = 4.16
$ws.Range("S6").Value = 1.77
$ws.Range("T6").Value = 0.07000000000000001
$ws.Range("V6").Value = 28.23
$ws.Range("W6").Value = 27.93
$ws.Range("X6").Value = 22.86
$ws.Range("Y6").Value = 0.32
$ws.Range("Z6").Value = 0.83
$ws.Range("AA6").Value = 4.58
$ws.Range("AB6").Value = 0.35
$ws.Range("AC6").Value = 0.07000000000000001
$ws.Range("AE6").Value = 34.88
$ws.Range("AF6").Value = 0.55
$ws.Range("AG6").Value = 1.63
$ws.Range("AH6").Value = 0.32
$ws.Range("AI6").Value = 0.58
$ws.Range("AJ6").Value = 0.32
$ws.Range("AL6").Value = 2.21

# Recomputed values for row 7 (Fixation duration (ms))
$ws.Range("D7").Value = 315.15
$ws.Range("F7").Value = 409.18
$ws.Range("G7").Value = 342.76
$ws.Range("H7").Value = 301.39
$ws.Range("I7").Value = 300.78
$ws.Range("J7").Value = 317.92
$ws.Range("K7").Value = 328.37
$ws.Range("M7").Value = 395.62
$ws.Range("N7").Value = 366.92
$ws.Range("Q7").Value = 516.1799999999999
$ws.Range("R7").Value = 358.16
$ws.Range("S7").Value = 326.54
$ws.Range("V7").Value = 300.38
$ws.Range("W7").Value = 303.4
$ws.Range("X7").Value = 316.08
$ws.Range("AA7").Value = 438.37
$ws.Range("AE7").Value = 409.97
$ws.Range("AG7").Value = 467.16

Write-Output "edit complete"
